# Weekly data refresh: insert two new "Acelga" price records at rows 693-694
# (date 2023-03-30), pushing all subsequent historical rows down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 693:694 - this shifts old rows 693..814 down to 695..816
$ws.Range("A693:A694").EntireRow.Insert()

# New row 693 - Primera quality record for 2023-03-30
$ws.Cells.Item(693, 1).Value = 9
$ws.Cells.Item(693, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(693, 3).Value = "Metropolitana"
$ws.Cells.Item(693, 4).Value = "2023-03-30"
$ws.Cells.Item(693, 5).Value = 13
$ws.Cells.Item(693, 6).Value = 100112009
$ws.Cells.Item(693, 7).Value = "Acelga"
$ws.Cells.Item(693, 8).Value = "Sin especificar"
$ws.Cells.Item(693, 9).Value = "Primera"
$ws.Cells.Item(693, 10).Value = 70
$ws.Cells.Item(693, 11).Value = 15000
$ws.Cells.Item(693, 12).Value = 15000
$ws.Cells.Item(693, 13).Value = 15000
$ws.Cells.Item(693, 14).Value = "`$/docena de atados"
$ws.Cells.Item(693, 15).Value = "Región Metropolitana"
$ws.Cells.Item(693, 16).Value = 5000
$ws.Cells.Item(693, 17).Value = 3
$ws.Cells.Item(693, 18).Value = "Hortaliza"

# New row 694 - Segunda quality record for 2023-03-30
$ws.Cells.Item(694, 1).Value = 9
$ws.Cells.Item(694, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(694, 3).Value = "Metropolitana"
$ws.Cells.Item(694, 4).Value = "2023-03-30"
$ws.Cells.Item(694, 5).Value = 13
$ws.Cells.Item(694, 6).Value = 100112009
$ws.Cells.Item(694, 7).Value = "Acelga"
$ws.Cells.Item(694, 8).Value = "Sin especificar"
$ws.Cells.Item(694, 9).Value = "Segunda"
$ws.Cells.Item(694, 10).Value = 34
$ws.Cells.Item(694, 11).Value = 13000
$ws.Cells.Item(694, 12).Value = 13000
$ws.Cells.Item(694, 13).Value = 13000
$ws.Cells.Item(694, 14).Value = "`$/docena de atados"
$ws.Cells.Item(694, 15).Value = "Región Metropolitana"
$ws.Cells.Item(694, 16).Value = 4333
$ws.Cells.Item(694, 17).Value = 3
$ws.Cells.Item(694, 18).Value = "Hortaliza"
